# Added Reference and Site entities, linked to Bionomics, with migration.
#
# This marks additional rows on the Bionomics sheet as already "Tabled"
# (column A = "y"), and introduces a new "Tabled?" marker column (header
# "Tabled", values "y") on the Site and Reference sheets, mirroring the
# column that already exists on the Bionomics sheet.

$wb = $excel.ActiveWorkbook

# --- Site sheet: add a "Tabled" column (A) with header + "y" for every data row ---
$wsSite = $wb.Worksheets.Item("Site")
$wsSite.Activate()
$wsSite.Range("A2").Value = "Tabled"
$wsSite.Range("A4:A25").Value = "y"
$wsSite.Range("A25").Select() | Out-Null

# --- Reference sheet: add a "Tabled" column (A) with header + "y" for every data row ---
$wsReference = $wb.Worksheets.Item("Reference")
$wsReference.Activate()
$wsReference.Range("A2").Value = "Tabled"
$wsReference.Range("A4:A11").Value = "y"
$wsReference.Range("A11").Select() | Out-Null

# --- Bionomics sheet: fill in two previously-missing "Tabled" markers ---
$wsBionomics = $wb.Worksheets.Item("Bionomics")
$wsBionomics.Activate()
$wsBionomics.Range("A5").Value = "y"
$wsBionomics.Range("A11").Value = "y"
$wsBionomics.Range("A12").Select() | Out-Null
